$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = '''2026-02-02'
$ws.Range("B2").Value = 'Organifarms'
$ws.Range("C2").Value = 'Robotics Software Engineer (m/w/d)'
$ws.Range("D2").Value = 'MATCH_SCORE: 8/10
GAPS: 1. Experience with Computer Vision Systems, 2. Knowledge of ROS packages such as MoveIt! and OMPL, 3. Practical experience with 7-Axis robots
ADVICE: To increase the chances of a successful application, the candidate should highlight any relevant projects or experiences that demonstrate their ability to learn and adapt to new technologies, such as ROS and Computer Vision Systems, and tailor their resume to emphasize transferable skills in robotics and software engineering.'
$ws.Range("E2").Value = 'https://www.linkedin.com/jobs/view/4357381390'

$ws.Range("A3").Value = '''2026-02-02'
$ws.Range("B3").Value = 'SiMa.ai'
$ws.Range("C3").Value = 'Principal, FAE (AI2472)'
$ws.Range("D3").Value = 'MATCH_SCORE: 6/10
GAPS: Technical knowledge of AI and ML, experience with Physical AI HW/SW platforms, familiarity with EMEA customer base
ADVICE: To increase the match score, the applicant should highlight any relevant experience or skills in AI, ML, or Physical AI platforms, and demonstrate an understanding of the EMEA market and customer needs.'
$ws.Range("E3").Value = 'https://www.linkedin.com/jobs/view/4368016450'

$ws.Range("A4").Value = '''2026-02-02'
$ws.Range("B4").Value = 'RoBoTec PTC'
$ws.Range("C4").Value = 'Modern C++ Developer (m/w/d)'
$ws.Range("D4").Value = 'MATCH_SCORE: 6/10
GAPS: 3D computer vision, robotics, high-performance system design 
ADVICE: To increase the chances of a successful application, the candidate should highlight any relevant experience or projects involving C++ and computer vision, and express a strong willingness to learn and adapt to the company''s specific technologies and challenges.'
$ws.Range("E4").Value = 'https://www.linkedin.com/jobs/view/4353892352'

$ws.Range("A5").Value = '''2026-02-02'
$ws.Range("B5").Value = 'SiMa.ai'
$ws.Range("C5").Value = 'ML Software Engineer (AI2464)'
$ws.Range("D5").Value = 'MATCH_SCORE: 6/10
GAPS: 1. Experience with Graph Neural Network algorithms, 2. Knowledge of SiMa SoC silicon and toolchain, 3. Fluency in German language
ADVICE: To increase the chances of a successful application, the candidate should highlight any transferable skills in ML engineering and AI research, and consider taking courses or certifications to address the identified skill gaps, particularly in Graph Neural Networks and SiMa SoC.'
$ws.Range("E5").Value = 'https://www.linkedin.com/jobs/view/4172879815'

$ws.Range("A6").Value = '''2026-02-02'
$ws.Range("B6").Value = 'SafeAD'
$ws.Range("C6").Value = 'Machine Learning Engineer'
$ws.Range("D6").Value = 'MATCH_SCORE: 8/10
GAPS: Experience with autonomous driving technology, Familiarity with Tensorflow or PyTorch for machine learning, Strong background in C++ programming
ADVICE: To increase the chances of a successful application, tailor the resume to highlight any relevant experience or projects related to machine learning, autonomous driving, or computer vision, and be prepared to explain how transferable skills can be applied to the role.'
$ws.Range("E6").Value = 'https://www.linkedin.com/jobs/view/4358376844'

$ws.Range("A7").Value = '''2026-02-02'
$ws.Range("B7").Value = 'SafeAD'
$ws.Range("C7").Value = 'Visual SLAM and 3D Reconstruction Engineer'
$ws.Range("D7").Value = 'MATCH_SCORE: 8/10
GAPS: Experience in autonomous driving, processing IMU and GPS data, and working on mid to large software projects
ADVICE: To increase the chances of getting hired, the applicant should highlight any relevant projects or experiences that demonstrate their ability to work with computer vision and machine learning concepts, and be prepared to discuss how their skills can be adapted to the field of autonomous driving.'
$ws.Range("E7").Value = 'https://www.linkedin.com/jobs/view/4358249370'

$ws.Range("A8").Value = '''2026-02-02'
$ws.Range("B8").Value = 'SafeAD'
$ws.Range("C8").Value = 'Visual SLAM and 3D Reconstruction Engineer'
$ws.Range("D8").Value = 'MATCH_SCORE: 60/100
GAPS: Experience in visual SLAM, sensor calibration, and autonomous driving, Strong programming skills in C++, and Experience working on mid to large software projects
ADVICE: To increase the chances of getting hired, the applicant should highlight any relevant projects or research experience in computer vision, machine learning, or robotics, and be prepared to demonstrate their problem-solving skills and learning mindset.'
$ws.Range("E8").Value = 'https://www.linkedin.com/jobs/view/4358249368'

$ws.Range("A9").Value = '''2026-02-02'
$ws.Range("B9").Value = 'ZEISS Group'
$ws.Range("C9").Value = 'Internship in Optical Sensing for Smart Instruments & Medical Robotics (f/m/x)'
$ws.Range("D9").Value = 'MATCH_SCORE: 8/10
GAPS: 1. Experience with optical coherence tomography (OCT), 2. Knowledge of medical robotics, 3. Familiarity with free-space and fiber-based optical systems
ADVICE: To increase the chances of a successful application, tailor your resume and cover letter to highlight any relevant coursework, projects, or research experience in optics, photonics, or biomedical engineering, and be prepared to discuss how your skills can be adapted to the specific requirements of the internship.'
$ws.Range("E9").Value = 'https://www.linkedin.com/jobs/view/4368156553'

$ws.Range("A10").Value = '''2026-02-02'
$ws.Range("B10").Value = 'Intrinsic'
$ws.Range("C10").Value = 'Intern: Open Source Developer Tools for Hard Real-Time Software Development'
$ws.Range("D10").Value = 'MATCH_SCORE: 6/10
GAPS: Experience building hard real-time systems in C++, familiarity with ROS2, and experience with build systems such as Bazel or CMake
ADVICE: To increase the chances of a successful application, tailor the resume to highlight any relevant coursework, personal projects, or experiences that demonstrate proficiency in C++ and interest in open-source software development, even if direct experience in hard real-time systems is limited.'
$ws.Range("E10").Value = 'https://www.linkedin.com/jobs/view/4357572585'

$ws.Range("A11").Value = '''2026-02-02'
$ws.Range("B11").Value = 'Innoviz Technologies'
$ws.Range("C11").Value = 'Tools and Automation Software Engineer'
$ws.Range("D11").Value = 'MATCH_SCORE: 60/100
GAPS: Experience with containerization technologies like Docker or Kubernetes, experience developing and maintaining automation pipelines using Jenkins and Groovy, and experience with robotics frameworks and middleware, including ROS (Robot Operating System)
ADVICE: To increase the chances of a successful application, the candidate should consider highlighting any transferable skills or experiences they have in related areas, such as automation or software development, and be prepared to address the gaps in their experience during the interview process.'
$ws.Range("E11").Value = 'https://www.linkedin.com/jobs/view/4346507333'

Write-Host "Done"